$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Step 1: add the new row (row 10) to the table, carrying over the row-above format ---
$lo.Resize($ws.Range("A1:E10"))
$ws.Range("A9:E9").Copy()
$ws.Range("A10:E10").PasteSpecial(-4122)

# Fill the new row 10 values (API reuses "main"; Pfad/Function/Body are brand-new strings)
$ws.Range("A10").Value2 = "main"
$ws.Range("B10").Value2 = "/review"
$ws.Range("D10").Value2 = "adds a new Review"
$ws.Range("E10").Value2 = "SecurityCookie, Product, Review"

# --- Step 2: add the new "request" column (F) to the table ---
$lo.Resize($ws.Range("A1:F10"))

# Header cell format (bold, like the other header cells)
$ws.Range("B1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data cells format for the new column (regular weight, no fill)
$ws.Range("B1").Copy()
$ws.Range("F2:F10").PasteSpecial(-4122)
$ws.Range("F2:F10").Font.Bold = $false

# Header text
$ws.Range("F1").Value2 = "request"

# Data values for the "main" API rows (5-10): get/get/get/get/put/get
$ws.Range("F5").Value2 = "get"
$ws.Range("F6").Value2 = "get"
$ws.Range("F7").Value2 = "get"
$ws.Range("F8").Value2 = "get"
$ws.Range("F9").Value2 = "put"
$ws.Range("F10").Value2 = "get"

# --- Step 3: update the Body column for the Bought/Delivered rows ---
$ws.Range("E8").Value2 = "SecurityCookie, Product"
$ws.Range("E9").Value2 = "SecurityCookie, Product"

# --- Step 4: widen column E to fit the new content ---
$ws.Range("E1").EntireColumn.ColumnWidth = 26.830729166666668

# --- Step 5: update the active selection ---
$ws.Range("E9").Select() | Out-Null
